$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.613.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.851.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.026"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.025"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4376"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3794"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07399"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.864.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.518"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.711"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07141"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "85.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.031"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009080"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.024"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.628.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.281"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.082.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.025"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.343"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.983"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08989"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7747"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.212"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  +3.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.554"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.026"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.139"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01972"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05263"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.856"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5181"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1679"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.834"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.835"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "110.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06597"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.027"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.701"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4700"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.892"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.56%  "
